# Insert a new data row at row 320 (pushing existing rows 320-424 down to
# 321-425) and populate it with the new "Cuatro cascos verde" record dated
# 2022-01-27 (serial 44588) for Región del Maule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("320:320").Insert()

$ws.Range("A320").Value = 5
$ws.Range("B320").Value = 'Macroferia Regional de Talca'
$ws.Range("C320").Value = 'Maule'
$ws.Range("D320").Value = 44588
$ws.Range("E320").Value = 7
$ws.Range("F320").Value = 100112002
$ws.Range("G320").Value = 'Pimiento'
$ws.Range("H320").Value = 'Cuatro cascos verde'
$ws.Range("I320").Value = 'Primera'
$ws.Range("J320").Value = 300
$ws.Range("K320").Value = 6000
$ws.Range("L320").Value = 6000
$ws.Range("M320").Value = 6000
$ws.Range("N320").Value = '$/caja 15 kilos'
$ws.Range("O320").Value = 'Región del Maule'
$ws.Range("P320").Value = 400
$ws.Range("Q320").Value = 15
$ws.Range("R320").Value = 'Hortaliza'
